$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 200 - this shifts the existing row 200 (and everything
# below it) down by one, matching the diff where old row N's data now lives
# in row N+1 for N = 200..263, and a brand-new weekly data point is written
# into the freshly inserted row 200.
$ws.Rows.Item(200).Insert()

# Populate the new row 200 with the new weekly record. Columns shared with
# the rest of the table (A,B,C,E,F,G,H,I,J,K,L,Q,R,T) keep the same
# boilerplate values as every other row for this market/product.
$ws.Range("A200").Value = 3
$ws.Range("B200").Value = "Femacal de La Calera"
$ws.Range("C200").Value = "Coquimbo"
$ws.Range("D200").Value = 44524
$ws.Range("E200").Value = 5
$ws.Range("F200").Value = "Fruta"
$ws.Range("G200").Value = 100108
$ws.Range("H200").Value = "Tropicales y subtropicales"
$ws.Range("I200").Value = 100108002
$ws.Range("J200").Value = "Mango"
$ws.Range("K200").Value = "Sin especificar"
$ws.Range("L200").Value = "Primera"
$ws.Range("M200").Value = 580
$ws.Range("N200").Value = 6000
$ws.Range("O200").Value = 6500
$ws.Range("P200").Value = 6276
$ws.Range("Q200").Value = "$/bandeja 4 kilos"
$ws.Range("R200").Value = "Perú"
$ws.Range("S200").Value = 1569
$ws.Range("T200").Value = 4
